$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 90: the timestamp in column A was re-aligned to 07:00:00 on the same date ---
$ws.Range("A90").Value = 45471.2916666667

# --- Append the newly scraped row (row 91) ---
$ws.Range("A91").Value = 45474.6130671296
# Match the date/time style used by the rest of column A (same cellXf as A90)
$ws.Range("A90").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$ws.Range("B91").Value = 4500
$ws.Range("C91").Value = 3.72000002861023
$ws.Range("D91").Value = 3.23000001907349
$ws.Range("E91").Value = 3.72000002861023
$ws.Range("F91").Value = 3.23000001907349

# G91 (adj_close) is stored as text in this workbook's convention (see G2, G3, ... which
# are all shared strings holding the numeric text). A leading apostrophe forces Excel to
# keep the numeric-looking entry as text; reapply the Normal style afterwards so the cell
# keeps the sheet's default (unstyled) look, like the other G-column cells.
$ws.Range("G91").Formula = "'3.23000001907349"
$ws.Range("G91").Style = "Normal"

$ws.Range("H91").Value = "ESPE.MI"
